$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new disc row 18: Metallica - Master of Puppets
$ws.Cells.Item(18, 1).Value = 1203
$ws.Cells.Item(18, 2).Value = " 'Master of Puppets'"
$ws.Cells.Item(18, 3).Value = " '2018-03-03'"
$ws.Cells.Item(18, 4).Value = " 'Digital'"
$ws.Cells.Item(18, 5).Value = " 'WAV'"
$ws.Cells.Item(18, 6).Value = " '4096 MB'"
$ws.Cells.Item(18, 7).Value = " ''"

# Add new disc row 20: AC/DC - Highway to Hell
$ws.Cells.Item(20, 1).Value = 8892
$ws.Cells.Item(20, 2).Value = " 'Highway to Hell'"
$ws.Cells.Item(20, 3).Value = " '1979-11-06'"
$ws.Cells.Item(20, 4).Value = " 'Physical'"
$ws.Cells.Item(20, 5).Value = " ''"
$ws.Cells.Item(20, 6).Value = " ''"
$ws.Cells.Item(20, 7).Value = " 'CD'"

# Update the selection to reflect where the user left off editing
$ws.Range("A22").Select()
